# Weekly update: insert the two newest records (row 513 & 514) at the top of
# the data block and push the previously-existing rows (old 513-531) down by
# two rows (they become rows 515-533).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 513:514 - this shifts existing rows 513-531
# down to 515-533, carrying their formatting (e.g. the date style on column D)
# along with them.
$ws.Rows("513:514").Insert()

# --- New row 513 ---
$ws.Range("A513").Value2 = 1
$ws.Range("B513").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C513").Value = "Arica y Parinacota"
$ws.Range("D513").Value2 = 45075
$ws.Range("E513").Value2 = 15
$ws.Range("F513").Value2 = 100114013
$ws.Range("G513").Value = "Zanahoria"
$ws.Range("H513").Value = "Sin especificar"
$ws.Range("I513").Value = "Primera"
$ws.Range("J513").Value2 = 40
$ws.Range("K513").Value2 = 14000
$ws.Range("L513").Value2 = 15000
$ws.Range("M513").Value2 = 14500
$ws.Range("N513").Value = "$/saco 25 kilos"
$ws.Range("O513").Value = "Provincia de Calama"
$ws.Range("P513").Value2 = 580
$ws.Range("Q513").Value2 = 25
$ws.Range("R513").Value = "Hortaliza"

# --- New row 514 ---
$ws.Range("A514").Value2 = 1
$ws.Range("B514").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C514").Value = "Arica y Parinacota"
$ws.Range("D514").Value2 = 45075
$ws.Range("E514").Value2 = 15
$ws.Range("F514").Value2 = 100114013
$ws.Range("G514").Value = "Zanahoria"
$ws.Range("H514").Value = "Sin especificar"
$ws.Range("I514").Value = "Primera"
$ws.Range("J514").Value2 = 50
$ws.Range("K514").Value2 = 14000
$ws.Range("L514").Value2 = 15000
$ws.Range("M514").Value2 = 14500
$ws.Range("N514").Value = "$/saco 25 kilos"
$ws.Range("O514").Value = "Valle de Camiña"
$ws.Range("P514").Value2 = 580
$ws.Range("Q514").Value2 = 25
$ws.Range("R514").Value = "Hortaliza"
